$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before F ("Locale"), shifting Parcel..Tax from F:N to G:O
$ws.Columns("F:F").Insert()

# New header cell for the inserted column
$ws.Range("F1").Value = "Locale"

# Row 2 now describes a different property; rewrite every cell explicitly
$ws.Range("A2").Value = "200 N Trig"
$ws.Range("B2").Value = "200 N TRIGG LLC"
$ws.Range("C2").Value = "126C-E-24.00--0"
$ws.Range("D2").Value = "200 TRIGG AVENUE NORTH"

# "70" looks numeric, force it to stay text (matches the rest of the sheet's
# inline-string cells), then strip the quote-prefix style back to default.
$ws.Range("E2").Value = "'70"
$ws.Range("E2").Style = "Normal"

# Locale .. Personal Property Value are blank for this row
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""

$ws.Range("K2").Value = "%"

$ws.Range("L2").Value = ""

$ws.Range("M2").Formula = "=H2+I2"
$ws.Range("N2").Formula = "=M2*(K2/100)"
$ws.Range("O2").Formula = "=N2*(L2/100)"
